# Auto-generated edit script: append rows 220-257 to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new data rows (220-257) ---
# Row 220
$ws.Range("A220").Value = 43921
$ws.Range("B220").Value = "Qatif"
$ws.Range("C220").Value = "Eastern province"
$ws.Range("D220").Formula = "=D219+E220"
$ws.Range("E220").Value = 7
$ws.Range("H220").Value = "Ash Sharqiyah"
# Row 221
$ws.Range("A221").Value = 43921
$ws.Range("B221").Value = "Mecca"
$ws.Range("C221").Value = "makkah"
$ws.Range("D221").Formula = "=D220+E221"
$ws.Range("E221").Value = 20
$ws.Range("H221").Value = "makkah"
# Row 222
$ws.Range("A222").Value = 43921
$ws.Range("B222").Value = "jeddah"
$ws.Range("C222").Value = "makkah"
$ws.Range("D222").Formula = "=D221+E222"
$ws.Range("E222").Value = 29
$ws.Range("H222").Value = "makkah"
# Row 223
$ws.Range("A223").Value = 43921
$ws.Range("B223").Value = "riyadh"
$ws.Range("C223").Value = "riyadh"
$ws.Range("D223").Formula = "=D222+E223"
$ws.Range("E223").Value = 33
$ws.Range("H223").Value = "ar Riyad"
# Row 224
$ws.Range("A224").Value = 43921
$ws.Range("B224").Value = "Medinah"
$ws.Range("C224").Value = "Medinah"
$ws.Range("D224").Formula = "=D223+E224"
$ws.Range("E224").Value = 3
$ws.Range("H224").Value = "Al Madinah"
# Row 225
$ws.Range("A225").Value = 43921
$ws.Range("B225").Value = "Dammam"
$ws.Range("C225").Value = "eastern province"
$ws.Range("D225").Formula = "=D224+E225"
$ws.Range("E225").Value = 3
$ws.Range("H225").Value = "Ash Sharqiyah"
# Row 226
$ws.Range("A226").Value = 43921
$ws.Range("B226").Value = "kobar"
$ws.Range("C226").Value = "eastern province"
$ws.Range("D226").Formula = "=D225+E226"
$ws.Range("E226").Value = 4
$ws.Range("H226").Value = "Ash Sharqiyah"
# Row 227
$ws.Range("A227").Value = 43921
$ws.Range("B227").Value = "Abha"
$ws.Range("C227").Value = "Asir"
$ws.Range("D227").Formula = "=D226+E227"
$ws.Range("E227").Value = 1
$ws.Range("H227").Value = "``Asir"
# Row 228
$ws.Range("A228").Value = 43921
$ws.Range("B228").Value = "ahsaa"
$ws.Range("C228").Value = "Eastern province"
$ws.Range("D228").Formula = "=D227+E228"
$ws.Range("E228").Value = 2
$ws.Range("H228").Value = "ash Sharqiyah"
# Row 229
$ws.Range("A229").Value = 43921
$ws.Range("B229").Value = "al badea"
$ws.Range("C229").Value = "Riyadh"
$ws.Range("D229").Formula = "=D228+E229"
$ws.Range("E229").Value = 1
$ws.Range("H229").Value = "Ar Riyad"
# Row 230
$ws.Range("A230").Value = 43921
$ws.Range("B230").Value = "Jizan"
$ws.Range("C230").Value = "jizan"
$ws.Range("D230").Formula = "=D229+E230"
$ws.Range("E230").Value = 2
$ws.Range("H230").Value = "Jizan"
# Row 231
$ws.Range("A231").Value = 43921
$ws.Range("B231").Value = "Dhahran"
$ws.Range("C231").Value = "Eastern province"
$ws.Range("D231").Formula = "=D230+E231"
$ws.Range("E231").Value = 2
$ws.Range("H231").Value = "ash Sharqiyah"
# Row 232
$ws.Range("A232").Value = 43921
$ws.Range("B232").Value = "Khafji"
$ws.Range("C232").Value = "Eastern province"
$ws.Range("D232").Formula = "=D231+E232"
$ws.Range("E232").Value = 1
$ws.Range("H232").Value = "Ash Sharqiyah"
# Row 233
$ws.Range("A233").Value = 43921
$ws.Range("B233").Value = "Ras Tanura"
$ws.Range("C233").Value = "Eastern province"
$ws.Range("D233").Formula = "=D232+E233"
$ws.Range("E233").Value = 1
$ws.Range("H233").Value = "Ash Sharqiyah"
# Row 234
$ws.Range("A234").Value = 43921
$ws.Range("B234").Value = "Khamis Mushait"
$ws.Range("C234").Value = "Asir"
$ws.Range("D234").Formula = "=D233+E234"
$ws.Range("E234").Value = 1
$ws.Range("H234").Value = "``Asir"
# Row 235
$ws.Range("A235").Value = 43922
$ws.Range("B235").Value = "Riyadh"
$ws.Range("C235").Value = "Riyadh"
$ws.Range("D235").Formula = "=D234+E235"
$ws.Range("E235").Value = 7
$ws.Range("H235").Value = "Ar Riyad"
# Row 236
$ws.Range("A236").Value = 43922
$ws.Range("B236").Value = "Mecca"
$ws.Range("C236").Value = "Makkah"
$ws.Range("D236").Formula = "=D235+E236"
$ws.Range("E236").Value = 55
$ws.Range("H236").Value = "Makkah"
# Row 237
$ws.Range("A237").Value = 43922
$ws.Range("B237").Value = "Medinah"
$ws.Range("C237").Value = "Medinah"
$ws.Range("D237").Formula = "=D236+E237"
$ws.Range("E237").Value = 78
$ws.Range("H237").Value = "Al Madinah"
# Row 238
$ws.Range("A238").Value = 43922
$ws.Range("B238").Value = "Ahsaa"
$ws.Range("C238").Value = "Eastern province"
$ws.Range("D238").Formula = "=D237+E238"
$ws.Range("E238").Value = 3
$ws.Range("H238").Value = "Ash Sharqiyah"
# Row 239
$ws.Range("A239").Value = 43922
$ws.Range("B239").Value = "Jeddah"
$ws.Range("C239").Value = "Makkah"
$ws.Range("D239").Formula = "=D238+E239"
$ws.Range("E239").Value = 3
$ws.Range("H239").Value = "Makkah"
# Row 240
$ws.Range("A240").Value = 43922
$ws.Range("B240").Value = "Qatif"
$ws.Range("C240").Value = "Eastern province"
$ws.Range("D240").Formula = "=D239+E240"
$ws.Range("E240").Value = 6
$ws.Range("H240").Value = "ASh Sharqiyah"
# Row 241
$ws.Range("A241").Value = 43922
$ws.Range("B241").Value = "Al Henakiyah"
$ws.Range("C241").Value = "Medinah"
$ws.Range("D241").Formula = "=D240+E241"
$ws.Range("E241").Value = 1
$ws.Range("H241").Value = "Al Madinah"
# Row 242
$ws.Range("A242").Value = 43922
$ws.Range("B242").Value = "altaif"
$ws.Range("C242").Value = "Makkah"
$ws.Range("D242").Formula = "=D241+E242"
$ws.Range("E242").Value = 2
$ws.Range("H242").Value = "Makkah"
# Row 243
$ws.Range("A243").Value = 43922
$ws.Range("B243").Value = "Tabuk"
$ws.Range("C243").Value = "Tabuk"
$ws.Range("D243").Formula = "=D242+E243"
$ws.Range("E243").Value = 2
$ws.Range("H243").Value = "Tabouk"
# Row 244
$ws.Range("A244").Value = 43923
$ws.Range("B244").Value = "Mecca"
$ws.Range("C244").Value = "Makkah"
$ws.Range("D244").Formula = "=D243+E244"
$ws.Range("E244").Value = 48
$ws.Range("H244").Value = "Makkah"
# Row 245
$ws.Range("A245").Value = 43923
$ws.Range("B245").Value = "Medinah"
$ws.Range("C245").Value = "Medinah"
$ws.Range("D245").Formula = "=D244+E245"
$ws.Range("E245").Value = 46
$ws.Range("H245").Value = "Al Madinah"
# Row 246
$ws.Range("A246").Value = 43923
$ws.Range("B246").Value = "Jeddah"
$ws.Range("C246").Value = "Makkah"
$ws.Range("D246").Formula = "=D245+E246"
$ws.Range("E246").Value = 30
$ws.Range("H246").Value = "Makkah"
# Row 247
$ws.Range("A247").Value = 43923
$ws.Range("B247").Value = "Khafji"
$ws.Range("C247").Value = "Eastern province"
$ws.Range("D247").Formula = "=D246+E247"
$ws.Range("E247").Value = 9
$ws.Range("H247").Value = "Ash Sharqiyah"
# Row 248
$ws.Range("A248").Value = 43923
$ws.Range("B248").Value = "Riyadh"
$ws.Range("C248").Value = "Riyadh"
$ws.Range("D248").Formula = "=D247+E248"
$ws.Range("E248").Value = 7
$ws.Range("H248").Value = "Ar Riyad"
# Row 249
$ws.Range("A249").Value = 43923
$ws.Range("B249").Value = "Khamis Mushait"
$ws.Range("C249").Value = "Asir"
$ws.Range("D249").Formula = "=D248+E249"
$ws.Range("E249").Value = 6
$ws.Range("H249").Value = "``Asir"
# Row 250
$ws.Range("A250").Value = 43923
$ws.Range("B250").Value = "Qatif"
$ws.Range("C250").Value = "Eastern province"
$ws.Range("D250").Formula = "=D249+E250"
$ws.Range("E250").Value = 5
$ws.Range("H250").Value = "Ash Sharqiyah"
# Row 251
$ws.Range("A251").Value = 43923
$ws.Range("B251").Value = "Dhahran"
$ws.Range("C251").Value = "Eastern province"
$ws.Range("D251").Formula = "=D250+E251"
$ws.Range("E251").Value = 4
$ws.Range("H251").Value = "Ash Sharqiyah"
# Row 252
$ws.Range("A252").Value = 43923
$ws.Range("B252").Value = "Dammam"
$ws.Range("C252").Value = "Eastern province"
$ws.Range("D252").Formula = "=D251+E252"
$ws.Range("E252").Value = 4
$ws.Range("H252").Value = "Ash Sharqiyah"
# Row 253
$ws.Range("A253").Value = 43923
$ws.Range("B253").Value = "Abha"
$ws.Range("C253").Value = "Asir"
$ws.Range("D253").Formula = "=D252+E253"
$ws.Range("E253").Value = 2
$ws.Range("H253").Value = "``Asir"
# Row 254
$ws.Range("A254").Value = 43923
$ws.Range("B254").Value = "Kobar"
$ws.Range("C254").Value = "Eastern province"
$ws.Range("D254").Formula = "=D253+E254"
$ws.Range("E254").Value = 1
$ws.Range("H254").Value = "Ash Sharqiyah"
# Row 255
$ws.Range("A255").Value = 43923
$ws.Range("B255").Value = "Ras Tanura"
$ws.Range("C255").Value = "Eastern province"
$ws.Range("D255").Formula = "=D254+E255"
$ws.Range("E255").Value = 1
$ws.Range("H255").Value = "ASh Sharqiyah"
# Row 256
$ws.Range("A256").Value = 43923
$ws.Range("B256").Value = "Ahad Rafidah"
$ws.Range("C256").Value = "Asir"
$ws.Range("D256").Formula = "=D255+E256"
$ws.Range("E256").Value = 1
$ws.Range("H256").Value = "``Asir"
# Row 257
$ws.Range("A257").Value = 43923
$ws.Range("B257").Value = "Bisha"
$ws.Range("C257").Value = "Asir"
$ws.Range("D257").Formula = "=D256+E257"
$ws.Range("E257").Value = 1
$ws.Range("H257").Value = "``Asir"

# --- Apply date number formatting (style index 3, matching existing date cells) ---
$ws.Range("A219").Copy() | Out-Null
$ws.Range("A220").PasteSpecial(-4122)
$ws.Range("A221").PasteSpecial(-4122)
$ws.Range("A222").PasteSpecial(-4122)
$ws.Range("A223").PasteSpecial(-4122)
$ws.Range("A224").PasteSpecial(-4122)
$ws.Range("A225").PasteSpecial(-4122)
$ws.Range("A226").PasteSpecial(-4122)
$ws.Range("A227").PasteSpecial(-4122)
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("A230").PasteSpecial(-4122)
$ws.Range("A231").PasteSpecial(-4122)
$ws.Range("A232").PasteSpecial(-4122)
$ws.Range("A233").PasteSpecial(-4122)
$ws.Range("A234").PasteSpecial(-4122)
$ws.Range("A235").PasteSpecial(-4122)
$ws.Range("A236").PasteSpecial(-4122)
$ws.Range("A237").PasteSpecial(-4122)
$ws.Range("A238").PasteSpecial(-4122)
$ws.Range("A239").PasteSpecial(-4122)
$ws.Range("A240").PasteSpecial(-4122)
$ws.Range("A241").PasteSpecial(-4122)
$ws.Range("A242").PasteSpecial(-4122)
$ws.Range("A243").PasteSpecial(-4122)
$ws.Range("A244").PasteSpecial(-4122)
$ws.Range("A245").PasteSpecial(-4122)
$ws.Range("A246").PasteSpecial(-4122)
$ws.Range("A247").PasteSpecial(-4122)
$ws.Range("A248").PasteSpecial(-4122)
$ws.Range("A249").PasteSpecial(-4122)
$ws.Range("A250").PasteSpecial(-4122)
$ws.Range("A251").PasteSpecial(-4122)
$ws.Range("A252").PasteSpecial(-4122)
$ws.Range("A253").PasteSpecial(-4122)
$ws.Range("A254").PasteSpecial(-4122)
$ws.Range("A255").PasteSpecial(-4122)
$ws.Range("A256").PasteSpecial(-4122)
$ws.Range("A257").PasteSpecial(-4122)

# H233 carries the same date-style formatting as column A in the source diff
$ws.Range("A219").Copy() | Out-Null
$ws.Range("H233").PasteSpecial(-4122)

# G234 is an empty placeholder cell that still carries the date-style formatting
$ws.Range("A219").Copy() | Out-Null
$ws.Range("G234").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Restore view/selection state to match the saved workbook ---
$ws.Range("B88").Select()
